# :construction: Spanish wines added
#
# Adds six new wine/grape-variety rows (Spanish & one Portuguese wine) to the
# "European design. & varieties" worksheet, right after the existing last
# row (57 - "Latium (White)" / "Frascati").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("European design. & varieties")
$ws.Activate()

# New rows appended at the bottom of the table (rows 58-63).
# Column A = Region / Appellation, Column B = Grape varieties.
# Values are written in the exact order needed so the underlying shared
# string table grows with the same ordering as the source edit (row 60's
# grape-variety string is entered before its region string).
$ws.Range("A58").Value = "Ribera-del-duero (Red), Spain"
$ws.Range("B58").Value = "Tempranillo, Cabernet Sauvignon"

$ws.Range("A59").Value = "Penedès, Spain"
$ws.Range("B59").Value = "Grenache, Syrah, Monastrell, Merlot, Pinot Noir, Cabernet Sauvignon, Samsó, Ull de Llebre"

$ws.Range("B60").Value = "Tempranillo, Grenache, Cabernet sauvignon, Merlot"
$ws.Range("A60").Value = "Navarra (Red), Spain"

$ws.Range("A61").Value = "Rueda (White), Spain"
$ws.Range("B61").Value = "Verdejo"

$ws.Range("A62").Value = "Rías Baixas, Spain"
$ws.Range("B62").Value = "Albariño"

$ws.Range("A63").Value = "Barca Velha (Red), Portugal"
$ws.Range("B63").Value = "Touriga Franca, Touriga Nacional,Tinta Roriz, Tinto Cão"

# Scroll the view down a bit and move the selection to just past the new
# data, mirroring where the editor's cursor ended up after typing the rows.
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A64:A66").Select()
